{"js": "// Update the 25 two-digit multiplication problems in the practice table.\n//\n// We read the table's current values as a 2D array first (a full\n// snapshot), build the new 2D array from that snapshot, and only then\n// write it back in one shot. Deciding every new value from the\n// pre-edit snapshot \u2014 rather than, say, a sequence of independent\n// find/replace passes \u2014 matters here because a few of the *new* values\n// equal some *other* cell's *old* value (e.g. \"26\u00d773=\" and \"77\u00d753=\" each\n// appear once as an old value and once as a new value in the table), so\n// anything that could re-match an already-updated cell would corrupt\n// the result.\nconst replacements = {\n  \"30\u00d754=\": \"88\u00d742=\",\n  \"95\u00d772=\": \"13\u00d739=\",\n  \"27\u00d748=\": \"57\u00d728=\",\n  \"95\u00d739=\": \"52\u00d758=\",\n  \"82\u00d725=\": \"43\u00d790=\",\n  \"31\u00d762=\": \"26\u00d773=\",\n  \"95\u00d758=\": \"62\u00d759=\",\n  \"37\u00d765=\": \"97\u00d771=\",\n  \"77\u00d753=\": \"40\u00d791=\",\n  \"82\u00d723=\": \"38\u00d745=\",\n  \"40\u00d772=\": \"98\u00d762=\",\n  \"80\u00d757=\": \"77\u00d745=\",\n  \"12\u00d794=\": \"84\u00d786=\",\n  \"64\u00d747=\": \"42\u00d728=\",\n  \"44\u00d782=\": \"47\u00d729=\",\n  \"87\u00d771=\": \"80\u00d747=\",\n  \"26\u00d773=\": \"29\u00d737=\",\n  \"78\u00d795=\": \"65\u00d758=\",\n  \"54\u00d791=\": \"99\u00d780=\",\n  \"31\u00d712=\": \"43\u00d786=\",\n  \"78\u00d757=\": \"85\u00d728=\",\n  \"93\u00d797=\": \"25\u00d791=\",\n  \"61\u00d795=\": \"91\u00d732=\",\n  \"68\u00d727=\": \"77\u00d753=\",\n  \"93\u00d790=\": \"19\u00d775=\",\n};\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (const table of tables.items) {\n  table.load(\"values\");\n}\nawait context.sync();\n\nfor (const table of tables.items) {\n  const values = table.values;\n  const newValues = values.map((row) =>\n    row.map((cellText) =>\n      Object.prototype.hasOwnProperty.call(replacements, cellText)\n        ? replacements[cellText]\n        : cellText\n    )\n  );\n  table.values = newValues;\n}\nawait context.sync();\n", "ps1": "# Update the 25 two-digit multiplication problems in the practice table.\n#\n# Each cell is addressed positionally (row, column), and we capture every\n# cell's original text before writing anything. That matters because a\n# couple of the *new* values happen to equal some *other* cell's *old*\n# value (e.g. \"26\u00d773=\" and \"77\u00d753=\" both appear as an old value in one\n# cell and a new value in another) \u2014 a plain sequential document-wide\n# Find/Replace could re-match an already-updated cell. Reading all the\n# old values up front and then writing by position avoids that hazard.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# (row, col, expectedOldText, newText) for every data cell, 1-based indices.\n$edits = @(\n  @(1, 1, \"30\u00d754=\", \"88\u00d742=\"),\n  @(1, 2, \"95\u00d772=\", \"13\u00d739=\"),\n  @(1, 3, \"27\u00d748=\", \"57\u00d728=\"),\n  @(1, 4, \"95\u00d739=\", \"52\u00d758=\"),\n  @(1, 5, \"82\u00d725=\", \"43\u00d790=\"),\n  @(5, 1, \"31\u00d762=\", \"26\u00d773=\"),\n  @(5, 2, \"95\u00d758=\", \"62\u00d759=\"),\n  @(5, 3, \"37\u00d765=\", \"97\u00d771=\"),\n  @(5, 4, \"77\u00d753=\", \"40\u00d791=\"),\n  @(5, 5, \"82\u00d723=\", \"38\u00d745=\"),\n  @(10, 1, \"40\u00d772=\", \"98\u00d762=\"),\n  @(10, 2, \"80\u00d757=\", \"77\u00d745=\"),\n  @(10, 3, \"12\u00d794=\", \"84\u00d786=\"),\n  @(10, 4, \"64\u00d747=\", \"42\u00d728=\"),\n  @(10, 5, \"44\u00d782=\", \"47\u00d729=\"),\n  @(15, 1, \"87\u00d771=\", \"80\u00d747=\"),\n  @(15, 2, \"26\u00d773=\", \"29\u00d737=\"),\n  @(15, 3, \"78\u00d795=\", \"65\u00d758=\"),\n  @(15, 4, \"54\u00d791=\", \"99\u00d780=\"),\n  @(15, 5, \"31\u00d712=\", \"43\u00d786=\"),\n  @(20, 1, \"78\u00d757=\", \"85\u00d728=\"),\n  @(20, 2, \"93\u00d797=\", \"25\u00d791=\"),\n  @(20, 3, \"61\u00d795=\", \"91\u00d732=\"),\n  @(20, 4, \"68\u00d727=\", \"77\u00d753=\"),\n  @(20, 5, \"93\u00d790=\", \"19\u00d775=\")\n)\n\n# Pass 1: grab every target cell and its current text (read-only).\n$cells = @()\nforeach ($edit in $edits) {\n  $row = $edit[0]\n  $col = $edit[1]\n  $cell = $t.Cell($row, $col)\n  $cells += , $cell\n}\n\n# Pass 2: write the new text into each cell using the snapshot above, so\n# no write can be mistaken for a not-yet-processed \"old\" value.\nfor ($i = 0; $i -lt $edits.Count; $i++) {\n  $newText = $edits[$i][3]\n  $cells[$i].Range.Text = $newText\n}\n"}
